# Rename the original sheet, add a second sheet with a similar
# "test data" table, and refresh the selections on both tabs.

$wb = $excel.ActiveWorkbook

# --- sheet 1: rename "test_sheet" -> "test_sheet_1" -----------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "test_sheet_1"

# --- sheet 2: new sheet "test_sheet_2" placed right after sheet 1 ---------
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "test_sheet_2"

$ws2.Range("A1").Value = "test_keyword"
$ws2.Range("B1").Value = "product_name"
$ws2.Range("C1").Value = "processor"
$ws2.Range("D1").Value = "ram"
$ws2.Range("E1").Value = "rom"

$ws2.Range("A2").Value = "test_1_2"
$ws2.Range("B2").Value = "Asus Vivobook 16x"
$ws2.Range("C2").Value = "AMD Ryzen 5600H"
$ws2.Range("D2").Value = "8 GB"
$ws2.Range("E2").Value = "500 GB"

$ws2.Range("A3").Value = "test_2_2"
$ws2.Range("B3").Value = "Dell Inspiron 3542"
$ws2.Range("C3").Value = "Intel Pentium Dual Core"
$ws2.Range("D3").Value = "8 GB"
$ws2.Range("E3").Value = "750 GB"

# Column widths to mirror the "best fit" look of sheet 1 (values are the
# COM ColumnWidth equivalents of the stored OOXML widths, which include a
# fixed ~5/6 character padding baked in by the host).
$ws2.Columns.Item(1).ColumnWidth = 11.276041666666666
$ws2.Columns.Item(2).ColumnWidth = 15.385416666666666
$ws2.Columns.Item(3).ColumnWidth = 19.608072916666668
$ws2.Columns.Item(4).ColumnWidth = 3.9440104166666665
$ws2.Columns.Item(5).ColumnWidth = 5.944010416666667

# --- sheet 1: selection becomes the whole data range, no longer tab-selected
$ws1.Range("A1:E3").Select() | Out-Null

# New sheet becomes the active/visible tab, with its own in-sheet
# selection resting one row below the data (A4).
$ws2.Activate() | Out-Null
$ws2.Range("A4").Select() | Out-Null
